$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.248.82"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "2.245.38"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.40%  "
$ws.Range("E6").Value = "  -3.84%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.527"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.51%  "
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.22"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.93%  "
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("D14").Value = "2.335.07"
$ws.Range("E14").Value = "  +3.91%  "
$ws.Range("D15").Value = "2.587.41"
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("E16").Value = "  -0.98%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.57"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.52%  "
$ws.Range("D18").Value = "44.005.79"
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("D19").Value = "0.0₃0965"
$ws.Range("E19").Value = "  -0.84%  "
$ws.Range("E20").Value = "  +1.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.63"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("E24").Value = "  -0.41%  "
$ws.Range("E25").Value = "  -0.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.10%  "
$ws.Range("E29").Value = "  +3.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.38%  "
$ws.Range("E31").Value = "  -1.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "153.52"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0798"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.14%  "
$ws.Range("E34").Value = "  -1.65%  "
$ws.Range("E35").Value = "  -4.14%  "
$ws.Range("E36").Value = "  +2.15%  "
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("E38").Value = "  -7.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.52"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.51%  "
$ws.Range("E40").Value = "  -4.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "14.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0301"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.20%  "
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("D44").Value = "1.750.32"
$ws.Range("E44").Value = "  +2.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "82.93"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.77%  "
$ws.Range("E46").Value = "  -1.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.95"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.20%  "
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("E50").Value = "  -1.91%  "
$ws.Range("E51").Value = "  -2.75%  "
